$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 51.05762969290213

$ws.Range("N2:N7").Value = $newValue
